# Cambio de formato para documento de aulas por defecto.
# Ahora, la primera columna sera el nombre del aula.
#
# The original sheet had headers: edificio | nombre | capacidad |
# equipamiento | horario_apertura | horario_cierre
#
# New layout: nombre | capacidad | equipamiento | horario_apertura |
# horario_cierre | edificio  (edificio moved from column A to column F)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "edificio " column (column A). The remaining headers
# shift left so column A now holds "nombre ", etc.
$ws.Columns.Item(1).Delete()

# Re-add the "edificio" header (no trailing space) as the new last column.
$ws.Range("F1").Value = "edificio"

# Touch A2 so a new (currently empty) second row exists, ready for data
# entry, and becomes the active selection - matching the new "type the
# classroom name first" workflow described in the commit message.
$ws.Range("A2").Style = "Normal"
$ws.Range("A2").Select()
